$wb = $excel.ActiveWorkbook

# ---- Sheet: VENTAS POR GRUPO ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L28").Value = 253.44
$ws1.Range("D30").Value = 292.99
$ws1.Range("D60").Value = "6 de 58"
$ws1.Range("L60").Value = "7 de 58"

# ---- Sheet: VENTA MENSUAL ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F28").Value = 3000.65
$ws2.Range("F30").Value = 681.79
$ws2.Range("F60").Value = 42520.12

# ---- Sheet: CUMPLIMIENTO MENSUAL ----
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 7114.66
$ws3.Range("E3").Value = 13272.8174217135
$ws3.Range("F3").Value = 0.3489720602914118

$ws3.Range("D11").Value = 12144.47
$ws3.Range("E11").Value = 7428.590249249699
$ws3.Range("F11").Value = 0.6204686362453484

$ws3.Range("D14").Value = 46534.06
$ws3.Range("E14").Value = 53363.93284188785
$ws3.Range("F14").Value = 0.4658157654243478
